# CCC19 Derived Variables Spreadsheet
# Add four new binary derived variables for VTE risk (X03a-X03d):
#   X03a VTE_risk_v2   - risk of VTE, alternate order of level definition
#   X03b VTE_risk_low  - low-risk VTE binary indicator
#   X03c VTE_risk_int  - intermediate-risk VTE binary indicator
#   X03d VTE_risk_high - high-risk VTE binary indicator
# These are inserted right after the existing "X03 VTE_risk" row (row 327),
# pushing the remaining rows (old X07..X11) down by four rows, and the
# backing Table1 / autofilter range is grown to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the four new rows directly under the existing X03 row.
$ws.Rows("328:331").Insert()

# --- Row 328: X03a / VTE_risk_v2 ---
$ws.Cells.Item(328, 1).Value2 = "X03a"
$ws.Cells.Item(328, 2).Value2 = "VTE_risk_v2"
$ws.Cells.Item(328, 3).Value2 = "Other"
$ws.Cells.Item(328, 4).Value2 = "Risk of VTE based on malignancy diagnosis, alternate order of level definition"
$ws.Cells.Item(328, 5).Value2 = "Low-risk VTE malignancy; High-risk VTE malignancy; Intermediate-risk VTE malignancy; Other heme malignancy; Other solid malignancy"

# --- Row 329: X03b / VTE_risk_low ---
$ws.Cells.Item(329, 1).Value2 = "X03b"
$ws.Cells.Item(329, 2).Value2 = "VTE_risk_low"
$ws.Cells.Item(329, 3).Value2 = "Other"
$ws.Cells.Item(329, 4).Value2 = "Patient has had at least one malignancy with low risk of VTE"
$ws.Cells.Item(329, 5).Value2 = "0 = No; 1 = Yes"

# --- Rows 330-331: X03c / X03d (variable # column filled first for both rows) ---
$ws.Cells.Item(330, 1).Value2 = "X03c"
$ws.Cells.Item(331, 1).Value2 = "X03d"

# variable-name column for both rows
$ws.Cells.Item(330, 2).Value2 = "VTE_risk_int"
$ws.Cells.Item(331, 2).Value2 = "VTE_risk_high"

# category column for both rows
$ws.Cells.Item(330, 3).Value2 = "Other"
$ws.Cells.Item(331, 3).Value2 = "Other"

# description column for both rows
$ws.Cells.Item(330, 4).Value2 = "Patient has had at least one malignancy with intermediate risk of VTE"
$ws.Cells.Item(331, 4).Value2 = "Patient has had at least one malignancy with high risk of VTE"

# values column for both rows
$ws.Cells.Item(330, 5).Value2 = "0 = No; 1 = Yes"
$ws.Cells.Item(331, 5).Value2 = "0 = No; 1 = Yes"

# Row heights: X03a's description wraps to 3 lines, X03b-X03d wrap to 1 line
# (matching the sibling rows' wrap heights already in the sheet).
$ws.Rows(328).RowHeight = 46
$ws.Rows(329).RowHeight = 16
$ws.Rows(330).RowHeight = 16
$ws.Rows(331).RowHeight = 16

# Grow Table1 (and its autofilter) to cover the four extra rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E337"))

# Leave the selection where the author left it after the edit.
$ws.Range("A331").Select()

Write-Output "Inserted X03a-X03d VTE risk derived variables; Table1 resized to A1:E337"
